$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2021.1666
$ws.Range("J112").Value = 2054.8215
$ws.Range("L112").Value = 6164.4645
$ws.Range("N112").Value = -8380.4645
$ws.Range("H118").Value = 568.2308
$ws.Range("I118").Value = 587.7778
$ws.Range("J118").Value = 524.25
$ws.Range("K118").Value = 1763.3334
$ws.Range("L118").Value = 1572.75
$ws.Range("M118").Value = -106.3334
$ws.Range("N118").Value = -4886.75
$ws.Range("H132").Value = 8068.913
$ws.Range("I132").Value = 9227.166999999999
$ws.Range("K132").Value = 27681.501
$ws.Range("M132").Value = -25151.501
$ws.Range("H137").Value = 1568519
$ws.Range("I137").Value = 3334352.2
$ws.Range("J137").Value = 10430.883
$ws.Range("K137").Value = 10003056.6
$ws.Range("L137").Value = 31292.649
$ws.Range("M137").Value = -10000506.6
$ws.Range("N137").Value = -36392.649
$ws.Range("H138").Value = 2754.4126
$ws.Range("J138").Value = 2249.681
$ws.Range("L138").Value = 6749.043
$ws.Range("N138").Value = -17029.043

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 3150
$ws.Range("I26").Value = 3150
$ws.Range("K26").Value = 3150
$ws.Range("M26").Value = -2820
$ws.Range("H32").Value = 3584.037
$ws.Range("I32").Value = 2999.8696
$ws.Range("J32").Value = 6943
$ws.Range("K32").Value = 2999.8696
$ws.Range("L32").Value = 6943
$ws.Range("M32").Value = -2712.8696
$ws.Range("N32").Value = -7517
$ws.Range("H63").Value = 1706.75
$ws.Range("J63").Value = 1701
$ws.Range("L63").Value = 1701
$ws.Range("N63").Value = -3073
$ws.Range("H66").Value = 1706.75
$ws.Range("J66").Value = 1701
$ws.Range("L66").Value = 8505
$ws.Range("N66").Value = -15369
$ws.Range("H74").Value = 255031.81
$ws.Range("I74").Value = 506440.9
$ws.Range("J74").Value = 3622.7273
$ws.Range("K74").Value = 506440.9
$ws.Range("L74").Value = 3622.7273
$ws.Range("M74").Value = -505566.9
$ws.Range("N74").Value = -5370.7273
$ws.Range("H77").Value = 255031.81
$ws.Range("I77").Value = 506440.9
$ws.Range("J77").Value = 3622.7273
$ws.Range("K77").Value = 2532204.5
$ws.Range("L77").Value = 18113.6365
$ws.Range("M77").Value = -2527836.5
$ws.Range("N77").Value = -26849.6365
$ws.Range("H80").Value = 84104.75
$ws.Range("J80").Value = 84104.75
$ws.Range("L80").Value = 84104.75
$ws.Range("N80").Value = -86100.75
$ws.Range("H83").Value = 84104.75
$ws.Range("J83").Value = 84104.75
$ws.Range("L83").Value = 252314.25
$ws.Range("N83").Value = -262298.25
$ws.Range("H102").Value = 6949.7
$ws.Range("I102").Value = 9299.4
$ws.Range("J102").Value = 4600
$ws.Range("K102").Value = 9299.4
$ws.Range("L102").Value = 4600
$ws.Range("M102").Value = -7677.4
$ws.Range("N102").Value = -7844
$ws.Range("H122").Value = 2854.8
$ws.Range("I122").Value = 3022.0557
$ws.Range("J122").Value = 1349.5
$ws.Range("K122").Value = 9066.167099999999
$ws.Range("L122").Value = 4048.5
$ws.Range("M122").Value = -6616.167099999999
$ws.Range("N122").Value = -8948.5
$ws.Range("H132").Value = 2003.7941
$ws.Range("I132").Value = 1312.6538
$ws.Range("K132").Value = 3937.9614
$ws.Range("M132").Value = -1407.9614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 37500
$ws.Range("J57").Value = 37500
$ws.Range("L57").Value = 37500
$ws.Range("N57").Value = -38940
$ws.Range("H82").Value = 43762.2
$ws.Range("I82").Value = 19879.445
$ws.Range("K82").Value = 19879.445
$ws.Range("M82").Value = -19496.445
$ws.Range("H85").Value = 43762.2
$ws.Range("I85").Value = 19879.445
$ws.Range("K85").Value = 19879.445
$ws.Range("M85").Value = -18553.445
$ws.Range("H105").Value = 28891382
$ws.Range("I105").Value = 2002189.8
$ws.Range("J105").Value = 62502870
$ws.Range("K105").Value = 2002189.8
$ws.Range("L105").Value = 62502870
$ws.Range("M105").Value = -2000442.8
$ws.Range("N105").Value = -62506364
$ws.Range("H133").Value = 24999.5
$ws.Range("J133").Value = 24999.5
$ws.Range("L133").Value = 24999.5
$ws.Range("N133").Value = -35119.5
$ws.Range("H134").Value = 3522.8857
$ws.Range("I134").Value = 3332.1785
$ws.Range("K134").Value = 9996.5355
$ws.Range("M134").Value = -7461.5355
$ws.Range("H136").Value = 37500
$ws.Range("J136").Value = 37500
$ws.Range("L136").Value = 37500
$ws.Range("N136").Value = -47700
$ws.Range("H138").Value = 60851.332
$ws.Range("J138").Value = 60851.332
$ws.Range("L138").Value = 60851.332
$ws.Range("N138").Value = -71131.33199999999
$ws.Range("H139").Value = 60587.168
$ws.Range("J139").Value = 60587.168
$ws.Range("L139").Value = 60587.168
$ws.Range("N139").Value = -70867.16800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1997.4
$ws.Range("I22").Value = 2296.75
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 2296.75
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -1946.75
$ws.Range("N22").Value = -1500
$ws.Range("H31").Value = 4886.222
$ws.Range("I31").Value = 3892.2222
$ws.Range("J31").Value = 5880.222
$ws.Range("K31").Value = 3892.2222
$ws.Range("L31").Value = 5880.222
$ws.Range("M31").Value = -3597.2222
$ws.Range("N31").Value = -6470.222
$ws.Range("H34").Value = 4886.222
$ws.Range("I34").Value = 3892.2222
$ws.Range("J34").Value = 5880.222
$ws.Range("K34").Value = 3892.2222
$ws.Range("L34").Value = 5880.222
$ws.Range("M34").Value = -3690.2222
$ws.Range("N34").Value = -6284.222
$ws.Range("H122").Value = 4637.2173
$ws.Range("I122").Value = 3825.0833
$ws.Range("J122").Value = 5523.1816
$ws.Range("K122").Value = 11475.2499
$ws.Range("L122").Value = 16569.5448
$ws.Range("M122").Value = -9025.249899999999
$ws.Range("N122").Value = -21469.5448
$ws.Range("H133").Value = 53999
$ws.Range("J133").Value = 53999
$ws.Range("L133").Value = 53999
$ws.Range("N133").Value = -59059

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2081.7334
$ws.Range("I132").Value = 1258.7778
$ws.Range("J132").Value = 3316.1667
$ws.Range("K132").Value = 11329.0002
$ws.Range("L132").Value = 29845.5003
$ws.Range("M132").Value = -8799.0002
$ws.Range("N132").Value = -34905.5003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3432.7
$ws.Range("J7").Value = 1490
$ws.Range("L7").Value = 1490
$ws.Range("N7").Value = -1714
$ws.Range("H16").Value = 943.7
$ws.Range("I16").Value = 943.7
$ws.Range("K16").Value = 943.7
$ws.Range("M16").Value = -773.7
$ws.Range("H40").Value = 4987.077
$ws.Range("I40").Value = 5065.6895
$ws.Range("K40").Value = 5065.6895
$ws.Range("M40").Value = -4929.6895
$ws.Range("H100").Value = 1714.2
$ws.Range("I100").Value = 1467.1666
$ws.Range("K100").Value = 1467.1666
$ws.Range("M100").Value = -926.1666
$ws.Range("H122").Value = 2575.2632
$ws.Range("I122").Value = 2698.5881
$ws.Range("J122").Value = 1527
$ws.Range("K122").Value = 8095.7643
$ws.Range("L122").Value = 4581
$ws.Range("M122").Value = -5645.7643
$ws.Range("N122").Value = -9481
$ws.Range("H126").Value = 3432.7
$ws.Range("J126").Value = 1490
$ws.Range("L126").Value = 4470
$ws.Range("N126").Value = -9410
$ws.Range("H132").Value = 5332.6
$ws.Range("I132").Value = 2160
$ws.Range("J132").Value = 6390.1333
$ws.Range("K132").Value = 6480
$ws.Range("L132").Value = 19170.3999
$ws.Range("M132").Value = -3950
$ws.Range("N132").Value = -24230.3999
$ws.Range("H134").Value = 103985
$ws.Range("J134").Value = 103985
$ws.Range("L134").Value = 103985
$ws.Range("N134").Value = -114125
$ws.Range("H135").Value = 35248.75
$ws.Range("J135").Value = 35248.75
$ws.Range("L135").Value = 35248.75
$ws.Range("N135").Value = -45388.75
$ws.Range("H136").Value = 6249.25
$ws.Range("I136").Value = 8498.75
$ws.Range("K136").Value = 25496.25
$ws.Range("M136").Value = -22946.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4133.1665
$ws.Range("J81").Value = 6299.5
$ws.Range("L81").Value = 12599
$ws.Range("N81").Value = -14721
$ws.Range("H84").Value = 4133.1665
$ws.Range("J84").Value = 6299.5
$ws.Range("L84").Value = 62995
$ws.Range("N84").Value = -73603
$ws.Range("H107").Value = 763.4737
$ws.Range("I107").Value = 635.2
$ws.Range("J107").Value = 1244.5
$ws.Range("K107").Value = 1905.6
$ws.Range("L107").Value = 3733.5
$ws.Range("M107").Value = 14.39999999999986
$ws.Range("N107").Value = -7573.5
$ws.Range("H122").Value = 25002468
$ws.Range("I122").Value = 2587.6667
$ws.Range("K122").Value = 7763.000100000001
$ws.Range("M122").Value = -5313.000100000001
$ws.Range("H126").Value = 2287.7144
$ws.Range("I126").Value = 1502.3334
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 4507.0002
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -2037.0002
$ws.Range("N126").Value = -25940

Write-Output "done"